$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 454; this shifts existing rows 454:549 down to 455:550
# and automatically extends the sheet dimension to A1:R550.
$ws.Rows("454:454").Insert()

# Fill in the new row 454 with the new data record (weekly price update).
$ws.Cells.Item(454, 1).Value = 4
$ws.Cells.Item(454, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(454, 3).Value = "Los Lagos"
$ws.Cells.Item(454, 4).Value = 45275
$ws.Cells.Item(454, 5).Value = 10
$ws.Cells.Item(454, 6).Value = 100112043
$ws.Cells.Item(454, 7).Value = "Pepino ensalada"
$ws.Cells.Item(454, 8).Value = "Sin especificar"
$ws.Cells.Item(454, 9).Value = "Primera"
$ws.Cells.Item(454, 10).Value = 400
$ws.Cells.Item(454, 11).Value = 18000
$ws.Cells.Item(454, 12).Value = 18000
$ws.Cells.Item(454, 13).Value = 18000
$ws.Cells.Item(454, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(454, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(454, 16).Value = 300
$ws.Cells.Item(454, 17).Value = 60
$ws.Cells.Item(454, 18).Value = "Hortaliza"
